$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.595.31"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.470.60"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9489"
$ws.Range("E5").Value = "  -5.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "281.73"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3716"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3199"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.95"
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.061"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06708"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.640"
$ws.Range("E13").Value = "  +4.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.34"
$ws.Range("E14").Value = "  +6.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.286"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.478.35"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001040"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05773"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.68"
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9611"
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.707"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.80"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.23"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.289"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.724.75"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.322"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.93"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.640.97"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.93"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.973"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.354"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8429"
$ws.Range("E33").Value = "  -6.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.652"
$ws.Range("E34").Value = "  +27.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07890"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06115"
$ws.Range("E36").Value = "  +7.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.942"
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.72"
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02073"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.125"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9610"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1905"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.374"
$ws.Range("E43").Value = "  -12.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5446"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.50"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.595"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.97"
$ws.Range("E47").Value = "  +11.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5361"
$ws.Range("E48").Value = "  +4.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.838"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06464"
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.054"
$ws.Range("E51").Value = "  +0.48%  "
